# Weekly update: a new price observation is inserted at row 13 (pushing
# every subsequent row down by one), matching the "Fruta / hortaliza,
# semanal" weekly-refresh pattern used across these sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13; all rows 13..43 shift to 14..44
# and existing formatting (e.g. the date style on column D) carries down.
$ws.Rows.Item(13).Insert()

# Populate the newly-inserted row 13 with the new weekly observation.
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44715
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100114007
$ws.Range("G13").Value = "Jengibre"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("N13").Value = "$/caja 13 kilos"
$ws.Range("O13").Value = "Perú"
$ws.Range("P13").Value = 1192
$ws.Range("Q13").Value = 13
$ws.Range("R13").Value = "Hortaliza"
